$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows' "runs","balls","fours","sixes" (columns C-F) values were
# rearranged between rows 3-8 and 10-13 (row 2 and row 9 stay the same).
# Build the new values per destination row, reading from the current
# (pre-edit) values of the corresponding source row, then write them
# back as text so the stored cell type remains a string.

$rowMap = @{
    3  = 6
    4  = 11
    5  = 10
    6  = 7
    7  = 4
    8  = 3
    10 = 8
    11 = 5
    12 = 13
    13 = 12
}

# Snapshot current values for columns C:F (runs, balls, fours, sixes) before
# overwriting anything, since several rows swap values with each other.
$snapshot = @{}
foreach ($r in 2..13) {
    $snapshot[$r] = @(
        $ws.Cells.Item($r, 3).Text,
        $ws.Cells.Item($r, 4).Text,
        $ws.Cells.Item($r, 5).Text,
        $ws.Cells.Item($r, 6).Text
    )
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    $cells = @(
        $ws.Cells.Item($destRow, 3),
        $ws.Cells.Item($destRow, 4),
        $ws.Cells.Item($destRow, 5),
        $ws.Cells.Item($destRow, 6)
    )
    for ($i = 0; $i -lt 4; $i++) {
        $cells[$i].NumberFormat = "@"
        $cells[$i].Value = [string]$vals[$i]
    }
}
